# Slide 7 ("Procedure Example - Parameters"): the "Content Placeholder 2"
# shape has a small Pascal code sample ending in:
#   ...
#      writeln(x);
#   end.
# Remove the redundant parentheses around the writeln argument, turning
# "writeln(x);" into "writeln x;" (the leading "(" becomes a plain space,
# the rest of the run - "x;" - is unchanged).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# There are two "(x);" substrings in this text box (the call "inc(x);" and
# the call "writeln(x);"). Skip past the first one so we land on the run
# that immediately follows "writeln".
$firstParen = $tr.Find("(x);")
$afterFirst = $firstParen.Start + $firstParen.Length
$writelnArg = $tr.Find("(x);", $afterFirst)

$writelnArg.Text = " x;"
